# Auto-generated edit script: updates cryptos price/volume data
# per commit 'Updated cryptos list on Mon Mar  6 18:58:04 UTC 2023 with GitHub Actions'
#
# Each target cell is forced to text type (NumberFormat '@') before the
# assignment so that numeric-looking strings (e.g. '289.44') are not
# silently converted to real numbers by Excel's Value setter, then
# ClearFormats() strips the now-unneeded '@' format so the cell keeps
# its original (default) style, matching the source workbook exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextCell "D2" "22.567.88"
Set-TextCell "D3" "1.578.08"
Set-TextCell "E3" "  +0.38%  "
Set-TextCell "E4" "  -0.02%  "
Set-TextCell "D6" "289.44"
Set-TextCell "E6" "  -0.52%  "
Set-TextCell "D7" "0.3728"
Set-TextCell "E7" "  +0.46%  "
Set-TextCell "D8" "48.43"
Set-TextCell "E8" "  -3.13%  "
Set-TextCell "D9" "0.3357"
Set-TextCell "E9" "  -0.58%  "
Set-TextCell "D10" "1.141"
Set-TextCell "E10" "  -0.65%  "
Set-TextCell "D11" "0.07514"
Set-TextCell "E11" "  -0.35%  "
Set-TextCell "E12" "  -0.01%  "
Set-TextCell "D13" "21.06"
Set-TextCell "E13" "  -0.57%  "
Set-TextCell "D14" "6.005"
Set-TextCell "E14" "  -0.32%  "
Set-TextCell "D15" "6.962"
Set-TextCell "E15" "  -0.07%  "
Set-TextCell "D16" "1.582.01"
Set-TextCell "E16" "  +0.65%  "
Set-TextCell "E17" "  +0.32%  "
Set-TextCell "D18" "88.70"
Set-TextCell "E18" "  -2.02%  "
Set-TextCell "D19" "0.06770"
Set-TextCell "E19" "  -0.03%  "
Set-TextCell "D20" "6.417"
Set-TextCell "E20" "  +1.15%  "
Set-TextCell "D22" "16.58"
Set-TextCell "E22" "  +0.94%  "
Set-TextCell "E23" "  -0.55%  "
Set-TextCell "D24" "22.557.21"
Set-TextCell "E24" "  +0.37%  "
Set-TextCell "E25" "  +1.58%  "
Set-TextCell "D26" "2.598"
Set-TextCell "E26" "  -0.77%  "
Set-TextCell "D27" "153.12"
Set-TextCell "E27" "  +2.72%  "
Set-TextCell "D28" "19.79"
Set-TextCell "E28" "  -1.19%  "
Set-TextCell "E29" "  -1.17%  "
Set-TextCell "D30" "124.47"
Set-TextCell "E30" "  -0.57%  "
Set-TextCell "D31" "1.757.75"
Set-TextCell "E31" "  +0.46%  "
Set-TextCell "D32" "1.057"
Set-TextCell "E32" "  -1.19%  "
Set-TextCell "D33" "6.205"
Set-TextCell "E33" "  +0.04%  "
Set-TextCell "E34" "  +0.16%  "
Set-TextCell "D35" "9.747"
Set-TextCell "E35" "  -0.53%  "
Set-TextCell "D36" "0.08328"
Set-TextCell "E36" "  -0.25%  "
Set-TextCell "D37" "0.02472"
Set-TextCell "E37" "  -0.21%  "
Set-TextCell "D38" "0.2291"
Set-TextCell "E38" "  -0.43%  "
Set-TextCell "D39" "5.433"
Set-TextCell "E39" "  -0.01%  "
Set-TextCell "D40" "0.06405"
Set-TextCell "E40" "  -2.12%  "
Set-TextCell "D41" "1.303"
Set-TextCell "E41" "  -4.25%  "
Set-TextCell "B42" "TheSandbox"
Set-TextCell "C42" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D42" "0.6358"
Set-TextCell "E42" "  +2.15%  "
Set-TextCell "B43" "Aptos"
Set-TextCell "C43" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D43" "11.38"
Set-TextCell "E43" "  +0.60%  "
Set-TextCell "E44" "  +0.02%  "
Set-TextCell "D45" "13.99"
Set-TextCell "E45" "  -0.88%  "
Set-TextCell "D46" "0.6174"
Set-TextCell "E46" "  +5.46%  "
Set-TextCell "E47" "  -0.32%  "
Set-TextCell "D48" "2.070"
Set-TextCell "E48" "  -0.07%  "
Set-TextCell "D49" "125.59"
Set-TextCell "E49" "  -2.75%  "
Set-TextCell "D50" "1.222"
Set-TextCell "E50" "  -0.54%  "
Set-TextCell "D51" "0.07288"
Set-TextCell "E51" "  -0.31%  "
